$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, copying the format from the neighboring
# header cell (G1) so it gets the same bold/border/alignment style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the corresponding data value in H2 (plain numeric, unstyled like F2/G2).
$ws.Range("H2").Value = 0
